$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: icon -> anime
$ws.Range("D1").Value = "anime"

# Category row: 圖標 -> 動畫
$ws.Range("D3").Value = "動畫"

# Slime rows (5,6,7): add "slime" tag in column D
$ws.Range("D5").Value = "slime"
$ws.Range("D6").Value = "slime"
$ws.Range("D7").Value = "slime"

# Bat rows (9,10,11): add "bat" tag in column D
$ws.Range("D9").Value = "bat"
$ws.Range("D10").Value = "bat"
$ws.Range("D11").Value = "bat"

# Update selection to match saved view state
$ws.Range("D11").Select()
